# Refresh the cryptocurrency table (Price / Volume(1h) columns) with the
# latest scrape, and fix the VeChain/EnergySwap row ordering (rows 44-45).
$ws = $excel.ActiveWorkbook.ActiveSheet

$updates = @(
    @{ Row = 2; D = '42.396.10'; E = '  +1.53%  ' },
    @{ Row = 3; D = '2.276.05'; E = '  +0.61%  ' },
    @{ Row = 4; D = '1.00'; E = '  +0.05%  ' },
    @{ Row = 5; D = '307.33'; E = '  +1.30%  ' },
    @{ Row = 6; D = '98.32'; E = '  +6.78%  ' },
    @{ Row = 7; E = '  -0.14%  ' },
    @{ Row = 8; E = '  +0.01%  ' },
    @{ Row = 9; E = '  +2.36%  ' },
    @{ Row = 10; D = '35.81'; E = '  +10.23%  ' },
    @{ Row = 11; E = '  +0.09%  ' },
    @{ Row = 12; E = '  -1.88%  ' },
    @{ Row = 13; E = '  +0.62%  ' },
    @{ Row = 14; D = '2.628.09'; E = '  +0.60%  ' },
    @{ Row = 15; E = '  +1.20%  ' },
    @{ Row = 16; D = '2.260.44'; E = '  -0.75%  ' },
    @{ Row = 17; D = '0.795'; E = '  +2.93%  ' },
    @{ Row = 18; D = '42.295.25'; E = '  +1.54%  ' },
    @{ Row = 19; D = '12.52'; E = '  +0.91%  ' },
    @{ Row = 20; E = '  +0.65%  ' },
    @{ Row = 21; D = '5.97'; E = '  +0.47%  ' },
    @{ Row = 22; E = '  +0.69%  ' },
    @{ Row = 23; D = '240.79'; E = '  +0.40%  ' },
    @{ Row = 24; E = '  +0.56%  ' },
    @{ Row = 25; D = '1.95'; E = '  +1.08%  ' },
    @{ Row = 26; E = '  -0.02%  ' },
    @{ Row = 27; E = '  -0.60%  ' },
    @{ Row = 28; D = '37.85'; E = '  +7.35%  ' },
    @{ Row = 29; E = '  -0.23%  ' },
    @{ Row = 30; E = '  +0.95%  ' },
    @{ Row = 31; D = '159.14'; E = '  -0.86%  ' },
    @{ Row = 32; E = '  -0.38%  ' },
    @{ Row = 33; E = '  +0.03%  ' },
    @{ Row = 34; D = '3.13'; E = '  +4.06%  ' },
    @{ Row = 35; D = '0.0741'; E = '  -0.30%  ' },
    @{ Row = 36; D = '16.99'; E = '  +0.01%  ' },
    @{ Row = 37; D = '2.38'; E = '  +0.70%  ' },
    @{ Row = 38; E = '  +1.05%  ' },
    @{ Row = 39; E = '  +2.87%  ' },
    @{ Row = 40; E = '  -1.22%  ' },
    @{ Row = 41; D = '4.11'; E = '  +5.09%  ' },
    @{ Row = 42; D = '2.41'; E = '  +13.89%  ' },
    @{ Row = 43; D = '1.997.03'; E = '  -0.75%  ' },
    @{ Row = 44; B = 'VeChain'; C = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D = '0.0286'; E = '  +2.48%  ' },
    @{ Row = 45; B = 'EnergySwap'; C = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D = '18.94'; E = '  -2.01%  ' },
    @{ Row = 46; E = '  +4.10%  ' },
    @{ Row = 47; E = '  -3.19%  ' },
    @{ Row = 48; D = '53.11'; E = '  +1.44%  ' },
    @{ Row = 49; E = '  +0.55%  ' },
    @{ Row = 50; D = '72.21'; E = '  +0.13%  ' },
    @{ Row = 51; D = '92.02'; E = '  +1.31%  ' }
)

foreach ($u in $updates) {
    foreach ($col in @("B", "C", "D", "E")) {
        if ($u.ContainsKey($col)) {
            $cellRef = "$col$($u.Row)"
            $val = $u[$col]
            # Columns in this sheet are plain text; force text storage (leading
            # apostrophe) when the literal looks like a bare number, so e.g.
            # "1.00" / "307.33" keep their printed form instead of becoming 1 / 307.33.
            if ($col -eq 'D' -and $val -match '^[+-]?[0-9]*\.?[0-9]+$') {
                $ws.Range($cellRef).Value = "'" + $val
            } else {
                $ws.Range($cellRef).Value = $val
            }
        }
    }
}